# Adds feature tickets UNI-0003 .. UNI-0009 to the "Main" bug-tracker sheet.
# This mirrors the commit "Added feature tickets : UNI-0003 to UNI-0009":
# seven new Feature rows are appended below the two existing rows (UNI-0001,
# UNI-0002), each Open/Unresolved, reported & assigned to rhdelaro, Low
# priority, with a ticket-specific summary/details pair. Everything else
# (Summary sheet COUNTIF totals, the five embedded charts, data validation,
# etc.) is formula/cache driven off this sheet and updates on recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Columns: A=Type  B=Key  C=Summary  D=Reporter  E=Assignee  F=Priority
#          G=Status  H=Resolution  I=Created  J=Updated  K=Notes
$rows = @(
    @{ Row=4;  Key="UNI-0003"; Summary="<I> Email/Feedback feature";
       Created="12/15/2013 09:23:49";
       Notes="User button for sending email/feedback" },
    @{ Row=5;  Key="UNI-0004"; Summary="<I> Saving/Loading/resetting progress";
       Created="12/15/2013 09:25:05";
       Notes="Feature to save/load/reset progress" },
    @{ Row=6;  Key="UNI-0005"; Summary="<I> Tick mark for solving levels";
       Created="12/15/2013 09:27:02";
       Notes="Display 'check' (completed) or 'star' (completed without mistake) for each of the completed levels" },
    @{ Row=7;  Key="UNI-0006"; Summary="<I> Achievements";
       Created="12/15/2013 09:30:43";
       Notes="Achievements/GameCenter(?)" },
    @{ Row=8;  Key="UNI-0007"; Summary="<I> Enable 'Settings' button";
       Created="12/15/2013 09:31:58";
       Notes="Setting sounds and/or brightness" },
    @{ Row=9;  Key="UNI-0008"; Summary="<I> Map of levels per difficulty";
       Created="12/15/2013 09:34:24";
       Notes="Levels map with completion tick marks" },
    @{ Row=10; Key="UNI-0009"; Summary="<I> Hints";
       Created="12/15/2013 09:35:36";
       Notes="Hints (for sale?)" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = "Feature"
    $ws.Range("B$n").Value = $r.Key
    $ws.Range("C$n").Value = $r.Summary
    $ws.Range("D$n").Value = "rhdelaro"
    $ws.Range("E$n").Value = "rhdelaro"
    $ws.Range("F$n").Value = "Low"
    $ws.Range("G$n").Value = "Open"
    $ws.Range("H$n").Value = "Unresolved"
    $ws.Range("I$n").Value = $r.Created
    $ws.Range("K$n").Value = $r.Notes
}

# Keep the sheet's active selection pointing at the next empty row, like the
# row just after the freshly entered data (K11), matching typical manual
# data-entry flow.
$ws.Range("K11").Select()
